$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string in cell A1
$ws.Range("A1").Value = "Datos actualizados a 21 de Agosto de 2020 a las 16:10"

# Update country statistics (refreshed data snapshot)
# Row 4
$ws.Range("B4").Value = 5747604
$ws.Range("C4").Value = 1332
$ws.Range("D4").Value = 3096369
$ws.Range("E4").Value = 2473778
$ws.Range("G4").Value = 33
$ws.Range("H4").Value = 177457

# Row 6
$ws.Range("B6").Value = 2925337
$ws.Range("C6").Value = 21008
$ws.Range("D6").Value = 2175492
$ws.Range("E6").Value = 694671
$ws.Range("G6").Value = 199
$ws.Range("H6").Value = 55174

# Row 22
$ws.Range("B22").Value = 231389
$ws.Range("C22").Value = 105
$ws.Range("E22").Value = 17264
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 9325

# Row 44
$ws.Range("A44").Value = "Emiratos Arabes Unidos"
$ws.Range("B44").Value = 66193
$ws.Range("C44").Value = 391
$ws.Range("D44").Value = 58296
$ws.Range("E44").Value = 7527
$ws.Range("G44").Value = 1
$ws.Range("H44").Value = 370

# Row 45
$ws.Range("A45").Value = "Guatemala"
$ws.Range("B45").Value = 65983
$ws.Range("D45").Value = 54351
$ws.Range("E45").Value = 9126
$ws.Range("H45").Value = 2506

# Row 50
$ws.Range("B50").Value = 55211
$ws.Range("C50").Value = 219
$ws.Range("D50").Value = 40473
$ws.Range("E50").Value = 12946
$ws.Range("G50").Value = 4
$ws.Range("H50").Value = 1792

# Row 64
$ws.Range("B64").Value = 34921
$ws.Range("C64").Value = 162
$ws.Range("D64").Value = 32682
$ws.Range("E64").Value = 1727
$ws.Range("G64").Value = 2
$ws.Range("H64").Value = 512

# Row 69
$ws.Range("B69").Value = 30378
$ws.Range("C69").Value = 169
$ws.Range("D69").Value = 28275
$ws.Range("E69").Value = 1411
$ws.Range("G69").Value = 3
$ws.Range("H69").Value = 692

# Row 76
$ws.Range("B76").Value = 18313
$ws.Range("C76").Value = 324
$ws.Range("D76").Value = 11102
$ws.Range("E76").Value = 7090
$ws.Range("G76").Value = 1
$ws.Range("H76").Value = 121

# Row 83
$ws.Range("B83").Value = 13308
$ws.Range("C83").Value = 114
$ws.Range("D83").Value = 9977
$ws.Range("E83").Value = 2774
$ws.Range("G83").Value = 3
$ws.Range("H83").Value = 557

# Row 96
$ws.Range("B96").Value = 8241
$ws.Range("C96").Value = 38
$ws.Range("D96").Value = 7040
$ws.Range("E96").Value = 1135

